# Fix naive component forecaster bug: insert a new first-lag error column
# before the existing ones (old B..J shift to C..K; values beyond K are dropped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.185867228879948
$ws.Range("C2").Value = -0.6149987959110895
$ws.Range("D2").Value = -1.573933492418202
$ws.Range("E2").Value = 1.259784732462786
$ws.Range("F2").Value = -0.4093199312433776
$ws.Range("G2").Value = 0.6108525064998651
$ws.Range("H2").Value = -1.102936950399446
$ws.Range("I2").Value = 1.725418423199526
$ws.Range("J2").Value = -0.5099141766665937
$ws.Range("K2").Value = 0.8414208812510687
$ws.Range("B3").Value = -2.762569267710482
$ws.Range("C3").Value = -3.721503964217595
$ws.Range("D3").Value = -0.8877857393366061
$ws.Range("E3").Value = -2.55689040304277
$ws.Range("F3").Value = -1.536717965299527
$ws.Range("G3").Value = -3.250507422198839
$ws.Range("H3").Value = -0.4221520485998669
$ws.Range("I3").Value = -2.657484648465986
$ws.Range("J3").Value = -1.306149590548324
$ws.Range("K3").Value = -1.927587520365226
$ws.Range("B4").Value = -0.281473977694265
$ws.Range("C4").Value = 2.552244247186724
$ws.Range("D4").Value = 0.8831395834805599
$ws.Range("E4").Value = 1.903312021223803
$ws.Range("F4").Value = 0.1895225643244911
$ws.Range("G4").Value = 3.017877937923463
$ws.Range("H4").Value = 0.7825453380573438
$ws.Range("I4").Value = 2.133880395975006
$ws.Range("J4").Value = 1.512442466158104
$ws.Range("K4").Value = 2.013003163348936
$ws.Range("B5").Value = 2.613418425600925
$ws.Range("C5").Value = 0.9443137618947609
$ws.Range("D5").Value = 1.964486199638004
$ws.Range("E5").Value = 0.2506967427386921
$ws.Range("F5").Value = 3.079052116337664
$ws.Range("G5").Value = 0.8437195164715449
$ws.Range("H5").Value = 2.195054574389207
$ws.Range("I5").Value = 1.573616644572305
$ws.Range("J5").Value = 2.074177341763137
$ws.Range("K5").Value = 1.406269269416995
$ws.Range("B6").Value = -1.097142175261494
$ws.Range("C6").Value = -0.07696973751825081
$ws.Range("D6").Value = -1.790759194417562
$ws.Range("E6").Value = 1.03759617918141
$ws.Range("F6").Value = -1.19773642068471
$ws.Range("G6").Value = 0.1535986372329528
$ws.Range("H6").Value = -0.467839292583949
$ws.Range("I6").Value = 0.032721404606882
$ws.Range("J6").Value = -0.6351866677392595
$ws.Range("K6").Value = -0.2061856925012563
$ws.Range("B7").Value = 0.5056247995153902
$ws.Range("C7").Value = -1.208164657383921
$ws.Range("D7").Value = 1.620190716215051
$ws.Range("E7").Value = -0.6151418836510686
$ws.Range("F7").Value = 0.7361931742665938
$ws.Range("G7").Value = 0.114755244449692
$ws.Range("H7").Value = 0.6153159416405229
$ws.Range("I7").Value = -0.05259213070561841
$ws.Range("J7").Value = 0.3764088445323847
$ws.Range("K7").Value = 0.4340803931105948
$ws.Range("B8").Value = -1.303839698193279
$ws.Range("C8").Value = 1.524515675405693
$ws.Range("D8").Value = -0.7108169244604263
$ws.Range("E8").Value = 0.6405181334572361
$ws.Range("F8").Value = 0.01908020364033419
$ws.Range("G8").Value = 0.5196409008311652
$ws.Range("H8").Value = -0.1482671715149762
$ws.Range("I8").Value = 0.2807338037230269
$ws.Range("J8").Value = 0.338405352301237
$ws.Range("K8").Value = 0.4623717231395225
$ws.Range("B9").Value = 1.784934712404416
$ws.Range("C9").Value = -0.4503978874617036
$ws.Range("D9").Value = 0.9009371704559588
$ws.Range("E9").Value = 0.279499240639057
$ws.Range("F9").Value = 0.780059937829888
$ws.Range("G9").Value = 0.1121518654837466
$ws.Range("H9").Value = 0.5411528407217497
$ws.Range("I9").Value = 0.5988243892999598
$ws.Range("J9").Value = 0.7227907601382453
$ws.Range("K9").Value = -0.1706273630965465
$ws.Range("B10").Value = -1.285852883620195
$ws.Range("C10").Value = 0.06548217429746761
$ws.Range("D10").Value = -0.5559557555194342
$ws.Range("E10").Value = -0.0553950583286032
$ws.Range("F10").Value = -0.7233031306747446
$ws.Range("G10").Value = -0.2943021554367415
$ws.Range("H10").Value = -0.2366306068585314
$ws.Range("I10").Value = -0.1126642360202459
$ws.Range("J10").Value = -1.006082359255038
$ws.Range("K10").Value = -0.2359516323112753
$ws.Range("B11").Value = 0.1655615342000891
$ws.Range("C11").Value = -0.4558763956168127
$ws.Range("D11").Value = 0.04468430157401831
$ws.Range("E11").Value = -0.6232237707721231
$ws.Range("F11").Value = -0.19422279553412
$ws.Range("G11").Value = -0.1365512469559099
$ws.Range("H11").Value = -0.01258487611762438
$ws.Range("I11").Value = -0.9060029993524162
$ws.Range("J11").Value = -0.1358722724086538
$ws.Range("K11").Value = -0.4017729932881683
$ws.Range("B12").Value = -0.1020898895371165
$ws.Range("C12").Value = 0.3984708076537146
$ws.Range("D12").Value = -0.2694372646924268
$ws.Range("E12").Value = 0.1595637105455762
$ws.Range("F12").Value = 0.2172352591237863
$ws.Range("G12").Value = 0.3412016299620719
$ws.Range("H12").Value = -0.55221649327272
$ws.Range("I12").Value = 0.2179142336710425
$ws.Range("J12").Value = -0.04798648720847212
$ws.Range("B13").Value = 0.2502934172212692
$ws.Range("C13").Value = -0.4176146551248722
$ws.Range("D13").Value = 0.0113863201131309
$ws.Range("E13").Value = 0.06905786869134101
$ws.Range("F13").Value = 0.1930242395296265
$ws.Range("G13").Value = -0.7003938837051653
$ws.Range("H13").Value = 0.06973684323859711
$ws.Range("I13").Value = -0.1961638776409175
$ws.Range("B14").Value = -0.5750606441290271
$ws.Range("C14").Value = -0.1460596688910241
$ws.Range("D14").Value = -0.08838812031281398
$ws.Range("E14").Value = 0.03557825052547153
$ws.Range("F14").Value = -0.8578398727093204
$ws.Range("G14").Value = -0.08770914576555788
$ws.Range("H14").Value = -0.3536098666450724
$ws.Range("B15").Value = 0.3545997876350467
$ws.Range("C15").Value = 0.4122713362132568
$ws.Range("D15").Value = 0.5362377070515423
$ws.Range("E15").Value = -0.3571804161832495
$ws.Range("F15").Value = 0.4129503107605129
$ws.Range("G15").Value = 0.1470495898809984
$ws.Range("B16").Value = 0.1319134556777877
$ws.Range("C16").Value = 0.2558798265160732
$ws.Range("D16").Value = -0.6375382967187186
$ws.Range("E16").Value = 0.1325924302250437
$ws.Range("F16").Value = -0.1333082906544708
$ws.Range("B17").Value = 0.4278546843610848
$ws.Range("C17").Value = -0.465563438873707
$ws.Range("D17").Value = 0.3045672880700554
$ws.Range("E17").Value = 0.03866656719054083
$ws.Range("B18").Value = -0.7714259786200386
$ws.Range("C18").Value = -0.001295251676276088
$ws.Range("D18").Value = -0.2671959725557906
$ws.Range("B19").Value = 0.6110347010110101
$ws.Range("C19").Value = 0.3451339801314955
$ws.Range("B20").Value = -0.343237405067616
